$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 9272.727999999999
$ws.Range("I21").Value = 2000
$ws.Range("J21").Value = 10000
$ws.Range("K21").Value = 2000
$ws.Range("L21").Value = 10000
$ws.Range("M21").Value = -1532
$ws.Range("N21").Value = -10936

$ws.Range("H23").Value = 9272.727999999999
$ws.Range("I23").Value = 2000
$ws.Range("J23").Value = 10000
$ws.Range("K23").Value = 2000
$ws.Range("L23").Value = 10000
$ws.Range("M23").Value = -1766
$ws.Range("N23").Value = -10468

$ws.Range("H41").Value = 430.13333
$ws.Range("J41").Value = 666.6667
$ws.Range("L41").Value = 666.6667
$ws.Range("N41").Value = -1546.6667

$ws.Range("H43").Value = 23811798
$ws.Range("I43").Value = 50002016
$ws.Range("J43").Value = 2508.4546
$ws.Range("K43").Value = 50002016
$ws.Range("L43").Value = 2508.4546
$ws.Range("M43").Value = -50001947
$ws.Range("N43").Value = -2646.4546

$ws.Range("H126").Value = 27499.25
$ws.Range("J126").Value = 27499.25
$ws.Range("L126").Value = 27499.25
$ws.Range("N126").Value = -37379.25

$ws.Range("H132").Value = 10876316
$ws.Range("I132").Value = 14292454
$ws.Range("J132").Value = 6785.8184
$ws.Range("K132").Value = 42877362
$ws.Range("L132").Value = 20357.4552
$ws.Range("M132").Value = -42874832
$ws.Range("N132").Value = -25417.4552

$ws.Range("H138").Value = 2638.8538
$ws.Range("I138").Value = 1423.862
$ws.Range("J138").Value = 3303.6604
$ws.Range("K138").Value = 4271.586
$ws.Range("L138").Value = 9910.9812
$ws.Range("M138").Value = 868.4139999999998
$ws.Range("N138").Value = -20190.9812

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12203730
$ws.Range("I32").Value = 15158544
$ws.Range("J32").Value = 15123.75
$ws.Range("K32").Value = 15158544
$ws.Range("L32").Value = 15123.75
$ws.Range("M32").Value = -15158257
$ws.Range("N32").Value = -15697.75

$ws.Range("H61").Value = 1606.4286
$ws.Range("I61").Value = 1431.025
$ws.Range("J61").Value = 1911.4783
$ws.Range("K61").Value = 1431.025
$ws.Range("L61").Value = 1911.4783
$ws.Range("M61").Value = -1219.025
$ws.Range("N61").Value = -2335.4783

$ws.Range("H74").Value = 15152685
$ws.Range("I74").Value = 16667886
$ws.Range("J74").Value = 671.3333
$ws.Range("K74").Value = 16667886
$ws.Range("L74").Value = 671.3333
$ws.Range("M74").Value = -16667012
$ws.Range("N74").Value = -2419.3333

$ws.Range("H77").Value = 15152685
$ws.Range("I77").Value = 16667886
$ws.Range("J77").Value = 671.3333
$ws.Range("K77").Value = 83339430
$ws.Range("L77").Value = 3356.6665
$ws.Range("M77").Value = -83335062
$ws.Range("N77").Value = -12092.6665

$ws.Range("H102").Value = 1115.4546
$ws.Range("I102").Value = 975.8421
$ws.Range("J102").Value = 1999.6666
$ws.Range("K102").Value = 975.8421
$ws.Range("L102").Value = 1999.6666
$ws.Range("M102").Value = 646.1579
$ws.Range("N102").Value = -5243.6666

$ws.Range("H122").Value = 1527.5
$ws.Range("I122").Value = 1496
$ws.Range("K122").Value = 4488
$ws.Range("M122").Value = -2038

$ws.Range("H136").Value = 1606.4286
$ws.Range("I136").Value = 1431.025
$ws.Range("J136").Value = 1911.4783
$ws.Range("K136").Value = 4293.075000000001
$ws.Range("L136").Value = 5734.4349
$ws.Range("M136").Value = -1743.075000000001
$ws.Range("N136").Value = -10834.4349

$ws.Range("H139").Value = 45675
$ws.Range("J139").Value = 45675
$ws.Range("L139").Value = 45675
$ws.Range("N139").Value = -55955

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1689732.9
$ws.Range("I22").Value = 1930923.2
$ws.Range("J22").Value = 1400
$ws.Range("K22").Value = 1930923.2
$ws.Range("L22").Value = 1400
$ws.Range("M22").Value = -1930750.2
$ws.Range("N22").Value = -1746

$ws.Range("H86").Value = 1059415.8
$ws.Range("I86").Value = 2911.25
$ws.Range("J86").Value = 2327221.2
$ws.Range("K86").Value = 2911.25
$ws.Range("L86").Value = 2327221.2
$ws.Range("M86").Value = -1788.25
$ws.Range("N86").Value = -2329467.2

$ws.Range("H89").Value = 1059415.8
$ws.Range("I89").Value = 2911.25
$ws.Range("J89").Value = 2327221.2
$ws.Range("K89").Value = 14556.25
$ws.Range("L89").Value = 11636106
$ws.Range("M89").Value = -8940.25
$ws.Range("N89").Value = -11647338

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1355.721
$ws.Range("I16").Value = 719.5
$ws.Range("J16").Value = 2159.3684
$ws.Range("K16").Value = 719.5
$ws.Range("L16").Value = 2159.3684
$ws.Range("M16").Value = -432.5
$ws.Range("N16").Value = -2733.3684

$ws.Range("H107").Value = 741.3333
$ws.Range("I107").Value = 608.5
$ws.Range("J107").Value = 807.75
$ws.Range("K107").Value = 608.5
$ws.Range("L107").Value = 807.75
$ws.Range("M107").Value = 1311.5
$ws.Range("N107").Value = -4647.75

$ws.Range("H113").Value = 1355.721
$ws.Range("I113").Value = 719.5
$ws.Range("J113").Value = 2159.3684
$ws.Range("K113").Value = 719.5
$ws.Range("L113").Value = 2159.3684
$ws.Range("M113").Value = 1450.5
$ws.Range("N113").Value = -6499.368399999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 697.92
$ws.Range("J131").Value = 776.525
$ws.Range("L131").Value = 2329.575
$ws.Range("N131").Value = -12409.575

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 19045.133
$ws.Range("J136").Value = 19045.133
$ws.Range("L136").Value = 57135.399
$ws.Range("N136").Value = -62235.399

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2555.5557

$ws.Range("H55").Value = 90918664
$ws.Range("I55").Value = 14792
$ws.Range("J55").Value = 250000450
$ws.Range("K55").Value = 14792
$ws.Range("L55").Value = 250000450
$ws.Range("M55").Value = -14619
$ws.Range("N55").Value = -250000796

$ws.Range("H126").Value = 2555.5557

$ws.Range("H132").Value = 16953622
$ws.Range("I132").Value = 32259928
$ws.Range("J132").Value = 7353.6787
$ws.Range("K132").Value = 96779784
$ws.Range("L132").Value = 22061.0361
$ws.Range("M132").Value = -96777254
$ws.Range("N132").Value = -27121.0361

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 17887.25
$ws.Range("J64").Value = 17887.25
$ws.Range("L64").Value = 17887.25
$ws.Range("N64").Value = -18383.25

$ws.Range("H67").Value = 17887.25
$ws.Range("J67").Value = 17887.25
$ws.Range("L67").Value = 17887.25
$ws.Range("N67").Value = -19603.25

$ws.Range("H96").Value = 2655.5557
$ws.Range("I96").Value = 2000
$ws.Range("J96").Value = 2842.8572
$ws.Range("K96").Value = 2000
$ws.Range("L96").Value = 2842.8572
$ws.Range("M96").Value = -627
$ws.Range("N96").Value = -5588.8572
